$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table by copying the existing 4 data rows (2-5) down into rows 6-9.
# This replicates formatting (e.g. the bold/bordered/centered style applied to
# column A) the same way Excel would when a user duplicates rows; afterwards we
# overwrite the values with the new game data.
$ws.Range("A2:AD5").Copy($ws.Range("A6"))

function Set-TextValue($cellAddr, $text) {
    # Assign as text (leading apostrophe prevents Excel from auto-converting
    # date-shaped strings like "2025-04-26" into a real date), then reset the
    # cell style back to Normal so no stray date number-format style lingers.
    $ws.Range($cellAddr).Value = "'" + $text
    $ws.Range($cellAddr).Style = "Normal"
}

# Row 6: HOU @ GSW (away), 2025-04-26
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "HOU"
$ws.Range("C6").Value = "GSW"
$ws.Range("D6").Value = "away"
Set-TextValue "E6" "2025-04-26"
$ws.Range("F6").Value = "240:00"
$ws.Range("G6").Value = 34
$ws.Range("H6").Value = 86
$ws.Range("I6").Value = 0.395
$ws.Range("J6").Value = 11
$ws.Range("K6").Value = 28
$ws.Range("L6").Value = 0.393
$ws.Range("M6").Value = 14
$ws.Range("N6").Value = 24
$ws.Range("O6").Value = 0.583
$ws.Range("P6").Value = 15
$ws.Range("Q6").Value = 36
$ws.Range("R6").Value = 51
$ws.Range("S6").Value = 19
$ws.Range("T6").Value = 5
$ws.Range("U6").Value = 6
$ws.Range("V6").Value = 14
$ws.Range("W6").Value = 17
$ws.Range("X6").Value = 93
$ws.Range("Y6").Value = -11
$ws.Range("Z6").Value = 22
$ws.Range("AA6").Value = 27
$ws.Range("AB6").Value = 22
$ws.Range("AC6").Value = 22
$ws.Range("AD6").Value = "L"

# Row 7: GSW vs HOU (home), 2025-04-26
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "GSW"
$ws.Range("C7").Value = "HOU"
$ws.Range("D7").Value = "home"
Set-TextValue "E7" "2025-04-26"
$ws.Range("F7").Value = "240:00"
$ws.Range("G7").Value = 39
$ws.Range("H7").Value = 90
$ws.Range("I7").Value = 0.433
$ws.Range("J7").Value = 14
$ws.Range("K7").Value = 45
$ws.Range("L7").Value = 0.311
$ws.Range("M7").Value = 12
$ws.Range("N7").Value = 15
$ws.Range("O7").Value = 0.8
$ws.Range("P7").Value = 10
$ws.Range("Q7").Value = 35
$ws.Range("R7").Value = 45
$ws.Range("S7").Value = 26
$ws.Range("T7").Value = 9
$ws.Range("U7").Value = 6
$ws.Range("V7").Value = 10
$ws.Range("W7").Value = 21
$ws.Range("X7").Value = 104
$ws.Range("Y7").Value = 11
$ws.Range("Z7").Value = 18
$ws.Range("AA7").Value = 28
$ws.Range("AB7").Value = 23
$ws.Range("AC7").Value = 35
$ws.Range("AD7").Value = "W"

# Row 8: HOU @ GSW (away), 2025-04-28
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "HOU"
$ws.Range("C8").Value = "GSW"
$ws.Range("D8").Value = "away"
Set-TextValue "E8" "2025-04-28"
$ws.Range("F8").Value = "240:00"
$ws.Range("G8").Value = 38
$ws.Range("H8").Value = 77
$ws.Range("I8").Value = 0.494
$ws.Range("J8").Value = 11
$ws.Range("K8").Value = 23
$ws.Range("L8").Value = 0.478
$ws.Range("M8").Value = 19
$ws.Range("N8").Value = 31
$ws.Range("O8").Value = 0.613
$ws.Range("P8").Value = 13
$ws.Range("Q8").Value = 28
$ws.Range("R8").Value = 41
$ws.Range("S8").Value = 18
$ws.Range("T8").Value = 6
$ws.Range("U8").Value = 6
$ws.Range("V8").Value = 13
$ws.Range("W8").Value = 21
$ws.Range("X8").Value = 106
$ws.Range("Y8").Value = -3
$ws.Range("Z8").Value = 26
$ws.Range("AA8").Value = 31
$ws.Range("AB8").Value = 23
$ws.Range("AC8").Value = 26
$ws.Range("AD8").Value = "L"

# Row 9: GSW vs HOU (home), 2025-04-28
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "GSW"
$ws.Range("C9").Value = "HOU"
$ws.Range("D9").Value = "home"
Set-TextValue "E9" "2025-04-28"
$ws.Range("F9").Value = "240:00"
$ws.Range("G9").Value = 36
$ws.Range("H9").Value = 86
$ws.Range("I9").Value = 0.419
$ws.Range("J9").Value = 17
$ws.Range("K9").Value = 46
$ws.Range("L9").Value = 0.37
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = 22
$ws.Range("O9").Value = 0.909
$ws.Range("P9").Value = 9
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = 34
$ws.Range("S9").Value = 27
$ws.Range("T9").Value = 7
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 8
$ws.Range("W9").Value = 22
$ws.Range("X9").Value = 109
$ws.Range("Y9").Value = 3
$ws.Range("Z9").Value = 28
$ws.Range("AA9").Value = 22
$ws.Range("AB9").Value = 32
$ws.Range("AC9").Value = 27
$ws.Range("AD9").Value = "W"
